$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value to a cell. Cells in this sheet store
# every field (including numbers) as plain text. Assigning a numeric-looking
# string straight to .Value would make Excel auto-convert it to a real
# number (and silently drop significant trailing zeros), so such cells are
# pre-formatted as Text ("@") to preserve the exact literal characters.
function Set-LiteralCellValue($Sheet, $CellRef, $Text) {
    $cell = $Sheet.Range($CellRef)
    if ($Text -match "^-?[0-9]+(\.[0-9]+)?$") {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $Text
}

$updates = @(
    @('D2', '244.38'),
    @('D3', '21.76'),
    @('D4', '5.390'),
    @('D5', '0.06013'),
    @('D6', '3.390'),
    @('D7', '0.8153'),
    @('D8', '0.9494'),
    @('B9', 'WazirX'),
    @('C9', 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'),
    @('D9', '0.1437'),
    @('E9', '8WazirXWRX'),
    @('B10', 'MandalaExchangeToken'),
    @('C10', 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'),
    @('D10', '0.07429'),
    @('E10', '9MandalaExchangeTokenMDX'),
    @('B11', 'LiechtensteinCryptoassetsExchange'),
    @('C11', 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'),
    @('D11', '0.03442'),
    @('E11', '10LiechtensteinCryptoassetsExchangeLCX'),
    @('B12', 'BitrueCoin'),
    @('C12', 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'),
    @('D12', '0.03056'),
    @('E12', '11BitrueCoinBTR'),
    @('B13', 'BitMartToken'),
    @('C13', 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'),
    @('D13', '0.09413'),
    @('E13', '12BitMartTokenBMX'),
    @('B14', 'MCDex'),
    @('C14', 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'),
    @('D14', '4.001'),
    @('E14', '13MCDexMCB'),
    @('B15', 'BitForexToken'),
    @('C15', 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'),
    @('D15', '0.001602'),
    @('E15', '14BitForexTokenBF'),
    @('B16', 'CoinExToken'),
    @('C16', 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'),
    @('D16', '0.04810'),
    @('E16', '15CoinExTokenCET'),
    @('B17', 'One'),
    @('C17', 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'),
    @('D17', '0.0005911'),
    @('E17', '16OneONE'),
    @('D18', '0.005597'),
    @('D19', '0.004161'),
    @('D20', '0.0009860'),
    @('D21', '3.666'),
    @('D22', '6.422'),
    @('D26', '0.00007001'),
    @('D40', '0.04008'),
    @('D41', '0.006441'),
    @('D42', '0.1073'),
    @('D44', '0.006659'),
    @('E44', '43LocalTradersLCTBestin24h'),
    @('D45', '0.00005236'),
    @('D48', '0.002975')
)

foreach ($u in $updates) {
    Set-LiteralCellValue $ws $u[0] $u[1]
}
